# EcoTaxisNYC_Sp3 - "PASOS A SEGUIR" slide: add MAE-related analysis bullets
# (commit: "agregado de analisis MAE")
#
# Slide 6 ("PASOS A SEGUIR" / "CONCLUSIONES"):
#   - Shape 3 ("CuadroTexto 5"): add two new numbered bullets right after the
#     first one, and grow/reposition the textbox to fit the extra lines.
#   - Shape 4 ("CONCLUSIONES" title textbox): move down to make room.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# ---------------------------------------------------------------------------
# 1) "CuadroTexto 5" - bulleted list of next steps
# ---------------------------------------------------------------------------
$bullets = $s.Shapes.Item(3)
$tf = $bullets.TextFrame
$tr = $tf.TextRange

# Insert the two new list items right after the first paragraph
# ("Entrenar al modelo con mas datos: ..."), pushing the remaining
# paragraphs ("Feature engineering: ...", etc.) further down.
$firstPara = $tr.Paragraphs(1, 1)
[void]$firstPara.InsertAfter("`rRealizar ajustes/mejoras en la arquitectura de la red neuronal.`rIncluir otros modelos: Sarimax, Prophet.")

# Split "Sarimax" and "Prophet" out into their own runs (they are flagged
# individually as product/model names in the source deck).
$trAfter = $tf.TextRange
$sarimaxRun = $trAfter.Find("Sarimax", 0)
$sarimaxRun.Text = "Sarimax"

$trAfter2 = $tf.TextRange
$prophetRun = $trAfter2.Find("Prophet", 0)
$prophetRun.Text = "Prophet"

# Resize/reposition the textbox so the two extra lines fit (grows upward
# and taller, same left/width). Values are in points (1 pt = 12700 EMU);
# chosen so the underlying EMU round-trips to the exact target offsets
# (x=1165214, y=649911, cx=9861572, cy=2862322).
$bullets.Left = 91.74917221069336
$bullets.Top = 51.17413330078125
$bullets.Width = 776.5017700195312
$bullets.Height = 225.37973022460938

# ---------------------------------------------------------------------------
# 2) "CONCLUSIONES" title textbox moves further down to stay clear of the
#    now-taller bullet list above (target y=3507122 EMU).
# ---------------------------------------------------------------------------
$conclusiones = $s.Shapes.Item(4)
$conclusiones.Top = 276.15138244628906
